# "Fixed typo on slides"
# Slide 33 ("Pattern Matching") has a code sample in "TextBox 2" that reads:
#   let result = divide 4 0
#   match x with
#   | None -> printfn "None"
#   | Some n -> printfn "Result: %i" n
#
# The "x" in "match x with" should actually refer to the "result" binding
# declared on the previous line, so fix the typo: " x " -> " result ".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(33)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$hit = $tr.Find(" x ")
if ($hit) {
    $hit.Text = " result "
    Write-Output "Slide 33: replaced ' x ' with ' result '."
} else {
    Write-Output "Slide 33: WARNING - ' x ' not found, no change made."
}
